$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G ("K") values are being regenerated for this save_data sheet.
# Update the affected rows to their new computed K values.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G6").Value = 1
